# ShadingPaints: Enlace entre la Vista de Pinturas y Datos de Pinturas
# Updates the paint-demand inputs (DEMANDA / HR_DISP / X / INV tables on Hoja1)
# and mirrors the user's scroll/selection position left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cellValues = @{
    "C5" = 500
    "E5" = 500
    "F5" = 500
    "G5" = 500
    "H5" = 500
    "I5" = 500
    "K5" = 500
    "L5" = 500
    "M5" = 500
    "N5" = 500
    "Q5" = 4
    "E6" = 500
    "F6" = 500
    "H6" = 500
    "I6" = 500
    "J6" = 500
    "K6" = 500
    "L6" = 500
    "N6" = 500
    "Q6" = 5
    "C7" = 500
    "D7" = 500
    "G7" = 500
    "H7" = 500
    "I7" = 500
    "J7" = 500
    "M7" = 500
    "N7" = 500
    "P7" = 3
    "Q7" = 4
    "C8" = 500
    "E8" = 500
    "G8" = 500
    "I8" = 500
    "J8" = 500
    "K8" = 500
    "L8" = 500
    "M8" = 500
    "Q8" = 3
    "C9" = 500
    "D9" = 500
    "E9" = 500
    "F9" = 500
    "G9" = 500
    "H9" = 500
    "I9" = 500
    "J9" = 500
    "K9" = 500
    "L9" = 500
    "M9" = 500
    "N9" = 500
    "P9" = 5
    "C10" = 500
    "D10" = 500
    "G10" = 500
    "H10" = 500
    "I10" = 500
    "M10" = 500
    "Q10" = 3
    "C11" = 500
    "D11" = 500
    "E11" = 500
    "F11" = 500
    "G11" = 500
    "H11" = 500
    "J11" = 500
    "K11" = 500
    "L11" = 500
    "M11" = 500
    "N11" = 500
    "Q11" = 5
    "C12" = 500
    "D12" = 500
    "F12" = 500
    "I12" = 500
    "K12" = 500
    "L12" = 500
    "N12" = 500
    "Q12" = 3
    "C13" = 500
    "E13" = 500
    "F13" = 500
    "G13" = 500
    "H13" = 500
    "I13" = 500
    "J13" = 500
    "K13" = 500
    "L13" = 500
    "M13" = 500
    "N13" = 500
    "P13" = 2
    "Q13" = 4
    "D14" = 500
    "E14" = 500
    "F14" = 500
    "K14" = 500
    "M14" = 500
    "N14" = 500
    "Q14" = 5
    "C15" = 500
    "E15" = 500
    "G15" = 500
    "H15" = 500
    "I15" = 500
    "J15" = 500
    "K15" = 500
    "L15" = 500
    "N15" = 500
    "P15" = 2
    "Q15" = 5
    "C16" = 500
    "D16" = 500
    "F16" = 500
    "G16" = 500
    "H16" = 500
    "I16" = 500
    "J16" = 500
    "K16" = 500
    "L16" = 500
    "M16" = 500
    "P16" = 1
    "Q16" = 4
    "C18" = 200000
    "D18" = 200000
    "E18" = 200000
    "F18" = 200000
    "G18" = 200000
    "H18" = 200000
    "I18" = 200000
    "J18" = 200000
    "K18" = 200000
    "L18" = 200000
    "M18" = 200000
    "N18" = 200000
    "C47" = 650
    "E47" = 500
    "F47" = 500
    "G47" = 500
    "H47" = 500
    "I47" = 500
    "K47" = 500
    "L47" = 500
    "M47" = 500
    "N47" = 500
    "D48" = 1000
    "F48" = 500
    "H48" = 500
    "I48" = 500
    "J48" = 500
    "K48" = 500
    "L48" = 500
    "N48" = 500
    "C49" = 1150
    "G49" = 500
    "H49" = 500
    "I49" = 500
    "J49" = 500
    "M49" = 500
    "N49" = 500
    "C50" = 650
    "E50" = 500
    "G50" = 500
    "I50" = 500
    "J50" = 500
    "K50" = 500
    "L50" = 500
    "M50" = 500
    "C51" = 650
    "D51" = 500
    "E51" = 2000
    "I51" = 500
    "J51" = 500
    "K51" = 500
    "L51" = 500
    "M51" = 500
    "N51" = 500
    "C52" = 650
    "D52" = 500
    "G52" = 500
    "H52" = 500
    "I52" = 500
    "M52" = 500
    "C53" = 650
    "D53" = 500
    "E53" = 500
    "F53" = 500
    "G53" = 500
    "H53" = 500
    "J53" = 500
    "K53" = 500
    "L53" = 500
    "M53" = 500
    "N53" = 500
    "C54" = 650
    "D54" = 500
    "F54" = 500
    "I54" = 500
    "K54" = 500
    "L54" = 500
    "N54" = 500
    "C55" = 650
    "E55" = 500
    "F55" = 500
    "G55" = 1000
    "I55" = 500
    "J55" = 500
    "K55" = 500
    "L55" = 500
    "M55" = 500
    "N55" = 500
    "D56" = 500
    "E56" = 500
    "F56" = 500
    "K56" = 500
    "M56" = 500
    "N56" = 500
    "C57" = 650
    "D57" = 1450
    "G57" = 500
    "H57" = 500
    "I57" = 500
    "J57" = 500
    "K57" = 500
    "L57" = 500
    "N57" = 500
    "C58" = 650
    "D58" = 500
    "F58" = 500
    "G58" = 500
    "H58" = 500
    "I58" = 500
    "J58" = 500
    "K58" = 500
    "L58" = 500
    "M58" = 500
    "D63" = 650
    "C64" = 650
    "E66" = 1650
    "F66" = 1150
    "G66" = 650
    "G70" = 650
    "D72" = 650
}

foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Value = $cellValues[$ref]
}

# Restore the view state (scroll position / active cell) left by the edit.
$ws.Activate()
$ws.Range("X67").Select()
